$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row needs to be inserted before the current row 787,
# shifting the existing rows 787-877 down to 788-878.
$ws.Rows("787:787").Insert()

$ws.Range("A787").Value = 10
$ws.Range("B787").Value = "Vega Modelo de Temuco"
$ws.Range("C787").Value = "La Araucanía"
$ws.Range("D787").Value = 45212
$ws.Range("E787").Value = 9
$ws.Range("F787").Value = "Fruta"
$ws.Range("G787").Value = 100101
$ws.Range("H787").Value = "Berries"
$ws.Range("I787").Value = 100101007
$ws.Range("J787").Value = "Kiwi"
$ws.Range("K787").Value = "Hayward"
$ws.Range("L787").Value = "Primera"
$ws.Range("M787").Value = 70
$ws.Range("N787").Value = 44000
$ws.Range("O787").Value = 44000
$ws.Range("P787").Value = 44000
$ws.Range("Q787").Value = "$/bandeja 18 kilos"
$ws.Range("R787").Value = "Región de O'Higgins"
$ws.Range("S787").Value = 2444
$ws.Range("T787").Value = 18
